# emily final chapter 8 edits
# Shift the dates in column F (rows 2-7) forward by one week (7 days),
# matching the original serial date values used in the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 44618
$ws.Range("F3").Value = 44617
$ws.Range("F4").Value = 44616
$ws.Range("F5").Value = 44615
$ws.Range("F6").Value = 44614
$ws.Range("F7").Value = 44613
